# Update workbook "phase_2_all_models_results" for Phase 2 results.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Key_Tenors_Comparison  (columns J = 5Y_pct, K = WAM_Years, rows 2-7)
# ---------------------------------------------------------------------------
$wsKey = $wb.Worksheets.Item("Key_Tenors_Comparison")

$wsKey.Range("J2").Value = [double]"1.11722290685571e-10"
$wsKey.Range("K2").Value = [double]"0.1319101748798887"

$wsKey.Range("J3").Value = [double]"54.50398119811268"
$wsKey.Range("K3").Value = [double]"2.689612997095686"

$wsKey.Range("J4").Value = [double]"45.76948907227623"
$wsKey.Range("K4").Value = [double]"2.59299930514257"

$wsKey.Range("J5").Value = [double]"57.62401814901956"
$wsKey.Range("K5").Value = [double]"2.806510561970689"

$wsKey.Range("J6").Value = [double]"51.77132656569123"
$wsKey.Range("K6").Value = [double]"2.519071808760003"

$wsKey.Range("J7").Value = [double]"69.83177860377482"
$wsKey.Range("K7").Value = [double]"3.643227801463834"

# ---------------------------------------------------------------------------
# Sheet: Full_Allocation  (columns F = Core_CF, H = Total_CF, J = CF_Percent)
# rows 12 (Exponential 5Y), 23 (Weibull 5Y), 34 (KaplanMeier 5Y),
# 45 (LogNormal 5Y), 56 (LogLogistic 5Y), 67 (Flat 5Y)
# ---------------------------------------------------------------------------
$wsFull = $wb.Worksheets.Item("Full_Allocation")

$wsFull.Range("F12").Value = [double]"2.083811187340613e-08"
$wsFull.Range("H12").Value = [double]"2.083811187340613e-08"
$wsFull.Range("J12").Value = [double]"1.11722290685571e-10"

$wsFull.Range("F23").Value = [double]"10165.92168655722"
$wsFull.Range("H23").Value = [double]"10165.92168655722"
$wsFull.Range("J23").Value = [double]"54.50398119811268"

$wsFull.Range("F34").Value = [double]"8536.789997986572"
$wsFull.Range("H34").Value = [double]"8536.789997986572"
$wsFull.Range("J34").Value = [double]"45.76948907227623"

$wsFull.Range("F45").Value = [double]"10747.8617688201"
$wsFull.Range("H45").Value = [double]"10747.8617688201"
$wsFull.Range("J45").Value = [double]"57.62401814901956"

$wsFull.Range("F56").Value = [double]"9656.235010851295"
$wsFull.Range("H56").Value = [double]"9656.235010851295"
$wsFull.Range("J56").Value = [double]"51.77132656569123"

$wsFull.Range("F67").Value = [double]"13024.81721360125"
$wsFull.Range("H67").Value = [double]"13024.81721360125"
$wsFull.Range("J67").Value = [double]"69.83177860377482"

# ---------------------------------------------------------------------------
# Sheet: 5Y_Ranking  (re-ranked by new Total_CF, descending)
# ---------------------------------------------------------------------------
$wsRank = $wb.Worksheets.Item("5Y_Ranking")

$wsRank.Range("A2").Value = "Flat"
$wsRank.Range("B2").Value = [double]"13024.81721360125"
$wsRank.Range("C2").Value = [double]"69.83177860377482"

$wsRank.Range("A3").Value = "LogNormal"
$wsRank.Range("B3").Value = [double]"10747.8617688201"
$wsRank.Range("C3").Value = [double]"57.62401814901956"

$wsRank.Range("A4").Value = "Weibull"
$wsRank.Range("B4").Value = [double]"10165.92168655722"
$wsRank.Range("C4").Value = [double]"54.50398119811268"

$wsRank.Range("A5").Value = "LogLogistic"
$wsRank.Range("B5").Value = [double]"9656.235010851295"
$wsRank.Range("C5").Value = [double]"51.77132656569123"

$wsRank.Range("A6").Value = "KaplanMeier"
$wsRank.Range("B6").Value = [double]"8536.789997986572"
$wsRank.Range("C6").Value = [double]"45.76948907227623"

$wsRank.Range("A7").Value = "Exponential"
$wsRank.Range("B7").Value = [double]"2.083811187340613e-08"
$wsRank.Range("C7").Value = [double]"1.11722290685571e-10"
